$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 44, shifting rows 44+ down by one.
$ws.Rows("44:44").Insert()

# Populate the newly inserted row 44 with the QUICKLOOK keyword data.
$ws.Range("A44").Value = "QUICKLOOK"
$ws.Range("B44").Value = $true
$ws.Range("C44").Value = "QUICKLOOK is true if you want to use the newer postprocessing system."

# Update the selection to match the saved state.
$ws.Range("C44").Select()
